# Apply the edits described in the commit: "Added images and 3 good templates"
# - Replace the header row / data block with a fresh 24-row template (new column A
#   header, and different X/Y/Angle minutiae values), trimming the old 32-row block
#   down to 24 rows.
# - Re-select the in-progress tagging cell (G11) as the active cell, matching the
#   saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:D1) ---------------------------------------------------
$ws.Range("A1").Value = " "
$ws.Range("B1").Value = "X"
$ws.Range("C1").Value = "Y"
$ws.Range("D1").Value = "Angle"

# --- Data rows (A2:D25) ----------------------------------------------------
$data = @(
  @(1, 153.878, 17.551, 25.56),
  @(2, 122.857, 49.796, 206.565),
  @(3, 172.245, 46.531, 17.103),
  @(4, 146.531, 69.388, 15.068),
  @(5, 181.224, 58.367, 198.435),
  @(6, 232.653, 85.306, 184.236),
  @(7, 235.102, 108.163, 2.49),
  @(8, 62.041, 181.633, 228.814),
  @(9, 186.122, 132.653, 187.431),
  @(10, 71.429, 215.918, 43.025),
  @(11, 105.306, 226.122, 40.914),
  @(12, 251.02, 170.612, 162.474),
  @(13, 50.612, 328.163, 43.152),
  @(14, 155.102, 298.776, 45),
  @(15, 74.694, 350.204, 37.694),
  @(16, 81.224, 362.857, 28.811),
  @(17, 121.633, 368.571, 33.69),
  @(18, 154.694, 345.714, 223.152),
  @(19, 133.469, 371.02, 28.811),
  @(20, 250.612, 342.449, 135),
  @(21, 274.286, 328.571, 117.897),
  @(22, 240.816, 355.918, 145.491),
  @(23, 199.184, 341.224, 20.854),
  @(24, 242.449, 389.388, 149.036)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# --- Drop the old trailing rows (26-33) that no longer exist ---------------
$ws.Range("A26:D33").ClearContents()

# --- Restore the saved selection/view state --------------------------------
$ws.Range("G11").Select() | Out-Null
